$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1840277777777778
$ws.Range("C2").Value = 0.5520833333333334
$ws.Range("J2").Value = 0.01736111111111111
$ws.Range("P2").Value = 0.1458333333333333
$ws.Range("S2").Value = 0.1006944444444444
$ws.Range("B3").Value = 0.01807228915662651
$ws.Range("C3").Value = 0.04819277108433735
$ws.Range("J3").Value = 0.03614457831325301
$ws.Range("P3").Value = 0.7168674698795181
$ws.Range("S3").Value = 0.1807228915662651
$ws.Range("J4").Value = 0.06666666666666667
$ws.Range("P4").Value = 0.7333333333333333
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.06896551724137931
$ws.Range("D6").Value = 0.009852216748768473
$ws.Range("F6").Value = 0.02955665024630542
$ws.Range("J6").Value = 0.2807881773399015
$ws.Range("O6").Value = 0.01970443349753695
$ws.Range("Q6").Value = 0.1428571428571428
$ws.Range("R6").Value = 0.03448275862068965
$ws.Range("S6").Value = 0.4137931034482759
$ws.Range("B7").Value = 0.1041666666666667
$ws.Range("D7").Value = 0.01041666666666667
$ws.Range("F7").Value = 0.03125
$ws.Range("J7").Value = 0.1302083333333333
$ws.Range("O7").Value = 0.05729166666666666
$ws.Range("Q7").Value = 0.1302083333333333
$ws.Range("R7").Value = 0.0625
$ws.Range("S7").Value = 0.4739583333333333
$ws.Range("B8").Value = 0.09843400447427293
$ws.Range("D8").Value = 0.01118568232662192
$ws.Range("F8").Value = 0.05369127516778523
$ws.Range("J8").Value = 0.1252796420581656
$ws.Range("O8").Value = 0.01565995525727069
$ws.Range("Q8").Value = 0.1633109619686801
$ws.Range("R8").Value = 0.07829977628635347
$ws.Range("S8").Value = 0.4541387024608501
$ws.Range("B9").Value = 0.1151515151515152
$ws.Range("D9").Value = 0.006060606060606061
$ws.Range("F9").Value = 0.04848484848484848
$ws.Range("J9").Value = 0.08484848484848485
$ws.Range("O9").Value = 0.0303030303030303
$ws.Range("Q9").Value = 0.1636363636363636
$ws.Range("R9").Value = 0.07272727272727272
$ws.Range("S9").Value = 0.4787878787878788
$ws.Range("B10").Value = 0.09951845906902086
$ws.Range("D10").Value = 0.01605136436597111
$ws.Range("E10").Value = 0.001605136436597111
$ws.Range("F10").Value = 0.06581059390048154
$ws.Range("J10").Value = 0.1252006420545747
$ws.Range("O10").Value = 0.008025682182985553
$ws.Range("Q10").Value = 0.2199036918138042
$ws.Range("R10").Value = 0.06179775280898876
$ws.Range("S10").Value = 0.4020866773675763
$ws.Range("G11").Value = 0.1661442006269593
$ws.Range("J11").Value = 0.1065830721003135
$ws.Range("K11").Value = 0.2068965517241379
$ws.Range("L11").Value = 0.5015673981191222
$ws.Range("S11").Value = 0.01880877742946709
$ws.Range("G12").Value = 0.7048192771084337
$ws.Range("J12").Value = 0.2168674698795181
$ws.Range("K12").Value = 0.01204819277108434
$ws.Range("L12").Value = 0.02409638554216868
$ws.Range("S12").Value = 0.04216867469879518
$ws.Range("G13").Value = 0.6
$ws.Range("J13").Value = 0.3142857142857143
$ws.Range("S13").Value = 0.08571428571428572
$ws.Range("F15").Value = 0.02325581395348837
$ws.Range("H15").Value = 0.1813953488372093
$ws.Range("I15").Value = 0.06046511627906977
$ws.Range("J15").Value = 0.3348837209302326
$ws.Range("K15").Value = 0.07906976744186046
$ws.Range("M15").Value = 0.004651162790697674
$ws.Range("O15").Value = 0.07906976744186046
$ws.Range("S15").Value = 0.2372093023255814
$ws.Range("F16").Value = 0.04444444444444445
$ws.Range("H16").Value = 0.1777777777777778
$ws.Range("I16").Value = 0.08888888888888889
$ws.Range("J16").Value = 0.3888888888888889
$ws.Range("M16").Value = 0.01111111111111111
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.1222222222222222
$ws.Range("F17").Value = 0.02147971360381861
$ws.Range("H17").Value = 0.1599045346062052
$ws.Range("I17").Value = 0.06682577565632458
$ws.Range("J17").Value = 0.4248210023866348
$ws.Range("K17").Value = 0.1026252983293556
$ws.Range("M17").Value = 0.01670644391408115
$ws.Range("O17").Value = 0.06443914081145585
$ws.Range("S17").Value = 0.1431980906921241
$ws.Range("F18").Value = 0.01398601398601399
$ws.Range("H18").Value = 0.1328671328671329
$ws.Range("I18").Value = 0.06293706293706294
$ws.Range("J18").Value = 0.4195804195804196
$ws.Range("K18").Value = 0.1258741258741259
$ws.Range("M18").Value = 0.01398601398601399
$ws.Range("O18").Value = 0.07692307692307693
$ws.Range("S18").Value = 0.1538461538461539
$ws.Range("F19").Value = 0.02323892519970951
$ws.Range("H19").Value = 0.2120551924473493
$ws.Range("I19").Value = 0.074800290486565
$ws.Range("J19").Value = 0.3558460421205519
$ws.Range("K19").Value = 0.1067538126361656
$ws.Range("M19").Value = 0.01815541031227306
$ws.Range("N19").Value = 0.002178649237472767
$ws.Range("O19").Value = 0.05954974582425563
$ws.Range("S19").Value = 0.1474219317356572
